# "Conquista de varias provincias a la vez y bucle de reseteo automatico"
#
# The "Control" sheet tracks, per province (rows 2-22), which team/club is
# currently in control (column B). This commit simulates a single team
# ("Minabo de Kiev") conquering every province at once, while column E keeps
# a log of who used to hold each group of provinces before the takeover.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Control")
$ws.Activate()

# New conqueror for every single province: column B becomes "Minabo de Kiev"
# for rows 2 through 22.
$ws.Range("B2:B22").Value = "Minabo de Kiev"

# Column E: history log of the previous controllers, grouped three rows at a
# time (the order the provinces were conquered in).
$previousControllers = @(
    "Al-Lagam",
    "Al-Lagam",
    "Al-Lagam",
    "Recreativo de Juerga",
    "Recreativo de Juerga",
    "Recreativo de Juerga",
    "Real Matriz",
    "Real Matriz",
    "Real Matriz",
    "Real Club de Parados",
    "Real Club de Parados",
    "Real Club de Parados",
    "Pombo FC",
    "Pombo FC",
    "Pombo FC",
    "Minabo de Kiev",
    "Minabo de Kiev",
    "Minabo de Kiev",
    "Gambote del Norte",
    "Gambote del Norte",
    "Gambote del Norte"
)

for ($i = 0; $i -lt $previousControllers.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $previousControllers[$i]
}

# Restore the selection/active-cell bookkeeping that Excel records after this
# bulk "select the whole result column" operation.
$ws.Range("B2:B22").Select()
